$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

$ws.Range("J12").Select()
